$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.598518443074256, 1.4010919640075317, -0.20246679710828802, 2.72812311031776)"
$ws.Range("C2").Value = "JSU(-1.6315669742205823, 1.5888618306996014, -0.4240514001540264, 8.340296966818615)"
$ws.Range("D2").Value = "JSU(-0.9356583481463515, 1.1299382344933107, 0.5188523317161331, 2.8022280446174896)"
$ws.Range("E2").Value = "NIG(0.7708529735243873, 0.5110552478964645, 5.233556553360591, 4.999476389601632)"
